# "Minor fixes on Dining room"
#  - Update the fixed "Date" footer placeholder text from 11/14/2020 to
#    11/20/2020 on the Slide Master and on every one of its 11 Custom
#    Layouts (the field itself is re-cached by PowerPoint whenever the
#    text is touched, which is the normal COM side-effect of editing a
#    header/footer date placeholder).
#  - On the room-map slide: the room that was still unlabeled ("?")
#    next to "Art Room" / "Arcade Room" is named "Music Hall", and the
#    room that used to be called "Music Hall" (next to "Main Hall" /
#    "Library") is renamed to "Longue".

$p = $ppt.ActivePresentation

# ---- 1. Fix the cached date footer text everywhere it appears ----
$newDate = "11/20/2020"

$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($l = 1; $l -le $m.CustomLayouts.Count; $l++) {
    $layout = $m.CustomLayouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---- 2. Rename the two room labels on slide 1 ----
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $t = $sh.TextFrame.TextRange.Text
        if ($sh.Name -eq "TextBox 12" -and $t -eq "?") {
            $sh.TextFrame.TextRange.Text = "Music Hall"
        } elseif ($sh.Name -eq "TextBox 58" -and $t -eq "Music Hall") {
            $sh.TextFrame.TextRange.Text = "Longue"
        }
    }
}
